$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "Save" header in H1, reusing the same formatting (bold,
# bordered, centered) as the other header cells by copying formats
# from the neighboring G1 header cell.
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Populate the new "Save" column for each data row (2-51) with the
# corresponding flag value.
$saveValues = @(0,0,0,0,0,0,0,0,0,0,0,0,0,1,0,1,0,0,0,0,1,0,0,0,0,0,0,0,0,1,0,1,0,0,0,0,0,0,0,1,0,0,0,0,0,0,0,0,0,0)

$startRow = 2
for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
